# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de sheets to reflect the newly generated report
# timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-24 10:23:45"
$zhcn.Range("H2").Value = "2016-03-24 10:24:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-24 10:23:50"
$dede.Range("H2").Value = "2016-03-24 10:24:19"
